# workers_rank_mat.xlsx: add infobox to make instructions for ranking
# easier/less (source commit re-shuffled a handful of worker rows plus
# refreshed every "mat_range" score). Columns: A=Unnamed:0 B=level_0
# C=index D=index(orig) E=prolificid F=name G=gender H=mat_range(score)
# I=race J=mat_rank K=matrices-range-bucket.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- matrices block (rows 2-13): three identities rotate between
#     rows 3/4/5 (index + prolificid + name), every row's score refreshed
$ws.Range("H2").Value  = 14.36475064273752

$ws.Range("D3").Value  = 2
$ws.Range("E3").Value  = "5f2c1a97a6809c060fec8820"
$ws.Range("F3").Value  = "Maggie"
$ws.Range("H3").Value  = 13.4427811560038

$ws.Range("D4").Value  = 1
$ws.Range("E4").Value  = "60b1742bce2b39e0f1d19a1a"
$ws.Range("F4").Value  = "Sabrina"
$ws.Range("H4").Value  = 13.32257368402617

$ws.Range("D5").Value  = 3
$ws.Range("E5").Value  = "60bd88b8fc436774352f53b9"
$ws.Range("F5").Value  = "Annes"
$ws.Range("H5").Value  = 13.02548504840682

$ws.Range("H6").Value  = 12.16366162123603
$ws.Range("H7").Value  = 10.35758251781631
$ws.Range("H8").Value  = 10.23661900101856
$ws.Range("H9").Value  = 8.201924197465678
$ws.Range("H10").Value = 5.441561929436489
$ws.Range("H11").Value = 2.330660576781288
$ws.Range("H12").Value = 2.005372734962068
$ws.Range("H13").Value = 1.34066941120993

# --- matrices block (rows 14-25): more identities rotate, genders follow
$ws.Range("H14").Value = 15.02328293437414
$ws.Range("H15").Value = 11.12005548300506

$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("F16").Value = "Alfredo"
$ws.Range("G16").Value = "male"
$ws.Range("H16").Value = 10.3560449567461

$ws.Range("D17").Value = 2
$ws.Range("E17").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("F17").Value = "Melissa"
$ws.Range("G17").Value = "female"
$ws.Range("H17").Value = 10.15590669353794

$ws.Range("H18").Value = 9.075645813370125
$ws.Range("H19").Value = 7.229575176107406

$ws.Range("D20").Value = 7
$ws.Range("E20").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("F20").Value = "Katherine"
$ws.Range("H20").Value = 5.144726965691964

$ws.Range("D21").Value = 9
$ws.Range("E21").Value = "5e35d91ea42bce592e996843"
$ws.Range("F21").Value = "Sergio"
$ws.Range("G21").Value = "male"
$ws.Range("H21").Value = 5.106254872490608

$ws.Range("D22").Value = 8
$ws.Range("E22").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("F22").Value = "Valeria"
$ws.Range("G22").Value = "female"
$ws.Range("H22").Value = 5.051234491524045

$ws.Range("H23").Value = 4.078136080597864
$ws.Range("H24").Value = 3.427904729701768
$ws.Range("H25").Value = 3.301880844181574
